# Insert a new weekly price-report row at row 244 (Ají, Vega Central Mapocho
# de Santiago), pushing the existing rows 244:300 down to 245:301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 244:300 down by one to make room for the new record.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A244").Value = 9
$ws.Range("B244").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C244").Value = "Metropolitana"
$ws.Range("D244").Value = 44722
$ws.Range("E244").Value = 13
$ws.Range("F244").Value = 100112021
$ws.Range("G244").Value = "Ají"
$ws.Range("H244").Value = "Americana (o)"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 16
$ws.Range("K244").Value = 25000
$ws.Range("L244").Value = 26000
$ws.Range("M244").Value = 25500
$ws.Range("N244").Value = "$/caja 25 kilos"
$ws.Range("O244").Value = "Provincia de Limarí"
$ws.Range("P244").Value = 1020
$ws.Range("Q244").Value = 25
$ws.Range("R244").Value = "Hortaliza"
